$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("250" -> "300"): update position size and rename tab ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C4").Value = 300
$ws1.Name = "300"
$ws1.Range("A1").Value = "比特币M5周期，每天1单，500点盈利"

# --- Sheet 2 ("1500"): fix the G4/G5:G18 divisor 600 -> 500 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G4").Formula = "=H4*E4/500"
$ws2.Range("G5:G18").Formula = "=H5*E5/500"
$ws2.Range("A1").Value = "比特币M5周期，每天1单，500点盈利"

# --- Sheet 3 ("5500"): fix the G4/G5:G18 divisor 600 -> 500 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("G4").Formula = "=H4*E4/500"
$ws3.Range("G5:G18").Formula = "=H5*E5/500"
$ws3.Range("A1").Value = "比特币M5周期，每天1单，500点盈利"

# --- Sheet 4 ("10000"): fix the G4/G5:G18 divisor 600 -> 500 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G4").Formula = "=H4*E4/500"
$ws4.Range("G5:G18").Formula = "=H5*E5/500"
$ws4.Range("A1").Value = "比特币M5周期，每天1单，500点盈利"

# --- Restore per-sheet selections (cursor position) ---
$ws1.Activate()
$ws1.Range("B4").Select()

$ws2.Activate()
$ws2.Range("B3").Select()

$ws3.Activate()
$ws3.Range("B23").Select()

# Sheet 4 ends up the active tab
$ws4.Activate()
$ws4.Range("B3").Select()
